$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Merge columns A and B into a single, wider column ---
$ws.Range("A1:B1").ColumnWidth = 34.45182291666666

# --- Add "v" marker column (V:W) for rows 2-9, reusing existing strings "CREAR"/"MARKUP" header in row 1 ---
$rows = 2..9
foreach ($r in $rows) {
    $ws.Cells.Item($r, 22).Value = "v"   # column V
    $ws.Cells.Item($r, 23).Value = "v"   # column W
}

# --- Spacer row ---
$ws.Range("A12").Font.Bold = $true

# --- Responsive footer height breakpoint table (rows 21-29) ---
$ws.Range("A21").Value = "GRAPHIC"
$ws.Range("B21").Value = 17
$ws.Range("H21").Value = 1280

$ws.Range("A22").Value = "T-SHIRT"
$ws.Range("B22").Value = 4
$ws.Range("C22").Value = 480
$ws.Range("D22").Value = 640

$ws.Range("A23").Value = "ILLUSTRATION"
$ws.Range("B23").Value = 8
$ws.Range("D23").Value = 640
$ws.Range("E23").Value = 768
$ws.Range("G23").Value = 1024
$ws.Range("H23").Value = 1280

$ws.Range("A24").Value = "POSTER"
$ws.Range("B24").Value = 12
$ws.Range("H24").Value = 1280

$ws.Range("A25").Value = "VISUAL IDENTITY"
$ws.Range("B25").Value = 1
$ws.Range("C25").Value = "x"

$ws.Range("A26").Value = "LETTERING"
$ws.Range("B26").Value = 6
$ws.Range("C26").Value = 480
$ws.Range("D26").Value = 640
$ws.Range("E26").Value = 768
$ws.Range("F26").Value = 800

$ws.Range("A27").Value = "COVER"
$ws.Range("B27").Value = 1
$ws.Range("C27").Value = "x"

$ws.Range("A28").Value = "SERIES"
$ws.Range("B28").Value = 2
$ws.Range("C28").Value = "x"

$ws.Range("A29").Value = "ARTWORK"
$ws.Range("B29").Value = 2
$ws.Range("C29").Value = "x"

# --- Section title row (added last so shared-string order matches) ---
$ws.Range("A19").Value = "RESPONSIVE FOOTER HEIGHT"
$ws.Range("A19").Font.Bold = $true

# --- Update selection / scroll position to match final state ---
$ws.Range("I21").Select() | Out-Null
